# "Update countries & provincias Spain"
# - Refresh a handful of per-country COVID stat rows (Albania, Montenegro,
#   Guyana, Butan).
# - Insert updated stats for "Republica de Africa Central" into the table,
#   which pushes Granada / Fiyi / Nueva Caledonia / Islas Virgenes de los
#   Estados Unidos / Namibia down by one row (Belice keeps its row/data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95: Albania
$ws.Range("A95").Value = "Albania"
$ws.Range("B95").Value = 726
$ws.Range("C95").Value = 14
$ws.Range("D95").Value = 410
$ws.Range("E95").Value = 288
$ws.Range("F95").Value = 4
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 28

# Row 121: Montenegro
$ws.Range("A121").Value = "Montenegro"
$ws.Range("B121").Value = 321
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 153
$ws.Range("E121").Value = 161
$ws.Range("F121").Value = 7
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 7

# Row 157: Guyana
$ws.Range("A157").Value = "Guyana"
$ws.Range("B157").Value = 74
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 12
$ws.Range("E157").Value = 54
$ws.Range("F157").Value = 5
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 8

# Row 183: Republica de Africa Central
$ws.Range("A183").Value = "Republica de Africa Central"
$ws.Range("B183").Value = 19
$ws.Range("C183").Value = 3
$ws.Range("D183").Value = 10
$ws.Range("E183").Value = 9
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

# Row 184: Belice
$ws.Range("A184").Value = "Belice"
$ws.Range("B184").Value = 18
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 5
$ws.Range("E184").Value = 11
$ws.Range("F184").Value = 1
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 2

# Row 185: Granada
$ws.Range("A185").Value = "Granada"
$ws.Range("B185").Value = 18
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 7
$ws.Range("E185").Value = 11
$ws.Range("F185").Value = 4
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0

# Row 186: Fiyi
$ws.Range("A186").Value = "Fiyi"
$ws.Range("B186").Value = 18
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 10
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# Row 187: Nueva Caledonia
$ws.Range("A187").Value = "Nueva Caledonia"
$ws.Range("B187").Value = 18
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 17
$ws.Range("E187").Value = 1
$ws.Range("F187").Value = 1
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

# Row 188: Islas Virgenes de los Estados Unidos
$ws.Range("A188").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B188").Value = 17
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 17
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

# Row 189: Namibia
$ws.Range("A189").Value = "Namibia"
$ws.Range("B189").Value = 16
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 8
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

# Row 190: Curazao
$ws.Range("A190").Value = "Curazao"
$ws.Range("B190").Value = 16
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 11
$ws.Range("E190").Value = 4
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 1

# Row 206: Butan
$ws.Range("A206").Value = "Butan"
$ws.Range("B206").Value = 7
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 4
$ws.Range("E206").Value = 3
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

Write-Output "done"